# CU1.4 Ver Venta.xlsx - apply commit "Modificacion descripcion, puntos de extension, curso basico"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0) Stash the cell formats (from rows 11-17) that do not have an equivalent
#    elsewhere in rows 1-10, because those rows are about to be deleted and
#    rebuilt. Stashing this way lets later copies reuse the existing style
#    table entries instead of Excel fabricating brand-new ones.
# ---------------------------------------------------------------------------
$ws.Range("A11:C11").Copy($ws.Range("A30:C30"))   # style (Curso Basico label row)
$ws.Range("A13:C13").Copy($ws.Range("A31:C31"))   # style (numbered wrapped-text step row)
$ws.Range("A15:C15").Copy($ws.Range("A32:C32"))   # style (blank spacer row, variant 1)
$ws.Range("A16:C16").Copy($ws.Range("A33:C33"))   # style (blank spacer row, variant 2)

# ---------------------------------------------------------------------------
# 1) Remove the old rows 11-17 entirely (this also drops any leftover custom
#    row heights so the rebuilt rows start out with the plain default height)
# ---------------------------------------------------------------------------
$ws.Range("11:17").Delete()

# ---------------------------------------------------------------------------
# 2) Rebuild rows 11-19 with the correct formatting, copied from rows that
#    already carry the right style (either still on the sheet or stashed).
# ---------------------------------------------------------------------------
$ws.Range("A9:C9").Copy($ws.Range("A11:C11"))      # Puntos de Extension value row
$ws.Range("A30:C30").Copy($ws.Range("A12:C12"))    # Curso Basico label row
$ws.Range("A31:C31").Copy($ws.Range("A13:C13"))    # Curso Basico step 1
$ws.Range("A31:C31").Copy($ws.Range("A14:C14"))    # Curso Basico step 2
$ws.Range("A31:C31").Copy($ws.Range("A15:C15"))    # Curso Basico step 3
$ws.Range("A4:C4").Copy($ws.Range("A16:C16"))      # Curso Alternativo label row
$ws.Range("A32:C32").Copy($ws.Range("A17:C17"))    # blank spacer row
$ws.Range("A33:C33").Copy($ws.Range("A18:C18"))    # blank spacer row
$ws.Range("A9:C9").Copy($ws.Range("A19:C19"))      # Pos-condicion value row

# ---------------------------------------------------------------------------
# 3) Clear the scratch/stash rows so they do not end up as stray content
# ---------------------------------------------------------------------------
$ws.Range("30:33").Delete()

# ---------------------------------------------------------------------------
# 4) Fill in the labels / values for rows 11-19
# ---------------------------------------------------------------------------

# Row 11 - Puntos de Extension: rich text with bold "labels"
$ws.Range("A11").ClearContents()
$cell = $ws.Range("B11")
$text = "Condicion: El actor quiere ver el detalle de una venta. Punto de extensión: Paso 2.d del CU01 Administracion de Supervisores: El actor hace clic en el botón ""Ver"" en algún registro del listado de ventas del formulario ""Administración de Supervisores""."
$cell.Value2 = $text

$run1 = $cell.Characters(1,10)
$run1.Font.Bold = $true
$run1.Font.Size = 10
$run1.Font.Name = "Arial"
$run1.Font.ColorIndex = -4105

$run2 = $cell.Characters(11,46)
$run2.Font.Bold = $false
$run2.Font.Size = 10
$run2.Font.Name = "Arial"
$run2.Font.ColorIndex = -4105

$run3 = $cell.Characters(57,18)
$run3.Font.Bold = $true
$run3.Font.Size = 10
$run3.Font.Name = "Arial"
$run3.Font.ColorIndex = -4105

$run4 = $cell.Characters(75,177)
$run4.Font.Bold = $false
$run4.Font.Size = 10
$run4.Font.Name = "Arial"
$run4.Font.ColorIndex = -4105

# Row 12 - Curso Basico label
$ws.Range("A12").Value2 = "Curso Básico"
$ws.Range("B12").ClearContents()

# Row 13 - Curso Basico step 1 (text unchanged, number goes from 2 to 1)
$ws.Range("A13").Value2 = 1
$ws.Range("B13").Value2 = "El sistema busca en la base de datos todos los datos coincidentes con el id de la venta seleccionada y los presenta en pantalla a traves del formulario Venta"

# Row 14 - Curso Basico step 2 (new)
$ws.Range("A14").Value2 = 2
$ws.Range("B14").Value2 = "El actor hace click en el botón ""Volver"""

# Row 15 - Curso Basico step 3 (new)
$ws.Range("A15").Value2 = 3
$ws.Range("B15").Value2 = "El sistema cierra el formulario."

# Row 16 - Curso Alternativo label
$ws.Range("A16").Value2 = "Curso Alternativo"
$ws.Range("B16").ClearContents()

# Row 17 / 18 - blank spacer rows, nothing else to set

# Row 19 - Pos-condicion
$ws.Range("A19").Value2 = "Pos-condición"
$ws.Range("B19").Value2 = "Se completaron todos los registros del formulario Venta."

# ---------------------------------------------------------------------------
# 5) Row height adjustments
# ---------------------------------------------------------------------------
$ws.Rows(7).RowHeight = 25.5
$ws.Rows(11).RowHeight = 51
$ws.Rows(13).RowHeight = 38.25

# ---------------------------------------------------------------------------
# 6) Simple value edits above the rebuilt block
# ---------------------------------------------------------------------------
$ws.Range("B5").Value2 = "0003"
$ws.Range("B7").Value2 = "Se muestran en pantalla todos los datos de la venta seleccionada incluyendo un botón ""Volver"""
$ws.Range("B9").Value2 = "Que el actor tenga los permisos necesarios para ver el registro." + [char]10 + "Que existan oportunidades con estado ""Ganada""."

# ---------------------------------------------------------------------------
# 7) View state: scroll + selection
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$ws.Range("B6").Select()
